$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.888.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.694.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +17.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "668.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  +6.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.55%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.692.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.20%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.385.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000268"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.648.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +17.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.705.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.537"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "516.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000209"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.97%  "
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.167"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.62%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.72%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.591"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "612.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.86%  "
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.963"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.88%  "
$ws.Range("E43").Value = "  +8.69%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0459"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.40%  "
$ws.Range("E47").Value = "  +26.42%  "
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.55%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.50%  "
